$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.971.16"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.625.32"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.30%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.49"
$ws.Range("E5").Value = "  -0.55%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -0.63%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.32%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.79%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -2.48%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.35"
$ws.Range("E10").Value = "  -5.36%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -0.75%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.852.86"
$ws.Range("E12").Value = "  -0.65%  "

# Rows 13 & 14 swap places: Polkadot (was row13) <-> WrappedEther (was row14)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.642.06"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.17"
$ws.Range("E14").Value = "  -1.51%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("E15").Value = "  -2.50%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.985.63"
$ws.Range("E16").Value = "  -0.01%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -1.86%  "

# Row 18 - Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.53"
$ws.Range("E18").Value = "  -2.51%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.29%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.04"
$ws.Range("E20").Value = "  -0.82%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.02%  "

# Row 22 - Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.53"
$ws.Range("E22").Value = "  -2.75%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.05"
$ws.Range("E23").Value = "  -1.66%  "

# Row 24 - Stellar
$ws.Range("E24").Value = "  +2.68%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.13"
$ws.Range("E25").Value = "  +0.52%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.40%  "

# Row 27 - Toncoin
$ws.Range("E27").Value = "  -2.66%  "

# Row 28 - Cosmos
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("E28").Value = "  -1.94%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.23"
$ws.Range("E29").Value = "  -1.44%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.32%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0483"
$ws.Range("E31").Value = "  -1.01%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.14"
$ws.Range("E32").Value = "  -2.73%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").Value = "  -4.11%  "

# Row 34 - HuobiToken
$ws.Range("E34").Value = "  -1.48%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -2.06%  "

# Row 36 - Maker
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.122.86"
$ws.Range("E36").Value = "  -0.40%  "

# Row 37 - ARBITRUM
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.851"
$ws.Range("E37").Value = "  -5.23%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  -1.11%  "

# Row 39 - ImmutableX
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.522"
$ws.Range("E39").Value = "  -1.60%  "

# Row 40 - VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0153"
$ws.Range("E40").Value = "  -1.64%  "

# Row 41 - Quant
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.25"
$ws.Range("E41").Value = "  -0.08%  "

# Row 42 - RocketPoolETH
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.762.47"
$ws.Range("E42").Value = "  -0.68%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.751"
$ws.Range("E43").Value = "  -5.12%  "

# Row 44 - FraxShare
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.10"
$ws.Range("E44").Value = "  -5.46%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  +1.37%  "

# Rows 46 & 47 swap places: Aave (was row46) <-> Cronos (was row47)
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0527"
$ws.Range("E46").Value = "  +1.29%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.31"
$ws.Range("E47").Value = "  -2.86%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.15%  "

# Row 49 - Mantle
$ws.Range("E49").Value = "  -0.16%  "

# Row 50 - EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.54"
$ws.Range("E50").Value = "  -1.54%  "

# Row 51 - USDD
$ws.Range("E51").Value = "  +0.42%  "
